$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "ppc" note for the Sparta row (both the 40% and 80% memory tables) ---
# The note text changes from "ppc=35" to "ppc=55"
$ws.Range("J14").Value = "ppc=55"
$ws.Range("J30").Value = "ppc=55"

# --- Update Sparta particles-per-cell derived value (reflects the new ppc=55 figure) ---
$ws.Range("C14").Value = 2522.1403676464702
$ws.Range("C30").Value = 2522.1403676464702

# --- Tidy up stray direct-formatting (font) that is not actually changing appearance ---
# These four cells only had a redundant "apply default font" flag; clear it.
$ws.Range("C11").ClearFormats()
$ws.Range("C14").ClearFormats()
$ws.Range("C27").ClearFormats()
$ws.Range("C30").ClearFormats()

# Re-apply the values (ClearFormats only touches formatting, but make sure content is intact)
$ws.Range("C14").Value = 2522.1403676464702
$ws.Range("C30").Value = 2522.1403676464702

# --- Update the saved selection / active cell for the sheet ---
$ws.Range("J31").Select()
